$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.846.86'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '1.624.34'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.84'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.29'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.74%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.256'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('D12').Value = '1.854.88'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').Value = '1.626.34'
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '27.845.35'
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.66'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('E30').Value = '  -1.33%  '
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('E33').Value = '  -0.51%  '
$ws.Range('D34').Value = '1.394.11'
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('E35').Value = '  +0.29%  '
$ws.Range('E36').Value = '  +11.03%  '
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  -1.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.849'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.56%  '
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('E44').Value = '  -2.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.77%  '
$ws.Range('D46').Value = '1.765.61'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('E47').Value = '  -2.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('E49').Value = '  +1.13%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0502'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.88%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.57'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.58%  '
